# Addes Start Date And End Date Of Available Days
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EnterShareSkill")

# H2/I2 held numeric Excel date-serials (11/25/2019) formatted with a date
# number format. They now hold explicit text dates for the skill's
# available-from / available-to range, so force the cells to Text format
# *before* assigning the strings (otherwise Excel auto-parses "10/12/2019"
# back into a date serial).
$ws.Range("H2").NumberFormat = "@"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("H2").Value = "10/12/2019"
$ws.Range("I2").Value = "15/12/2019"

# EnterShareSkill becomes the active sheet/tab (previously EditShareSkill),
# with a new selection.
$ws.Activate()
$ws.Range("F10:F11").Select()

# Page setup for EnterShareSkill.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
